# Applies the changes described in the commit:
#  - Switches the active sheet from "TimeIDDependentTime" to "EVData"
#  - Converts TimeIDDependentTime!B2:C25 from time-of-day fractions to full
#    datetime values (with a custom yyyy-mm-dd hh:mm:ss number format) and
#    merges/re-sizes columns B:C
#  - Updates the selections on both worksheets
#  - Highlights EVData!D1:D21 (MaxFastChargingPower) with a yellow fill

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# TimeIDDependentTime sheet: B2:C25 datetime values + formatting
# ---------------------------------------------------------------------
$wsTime = $wb.Worksheets.Item("TimeIDDependentTime")

$startSerial = 44348.291666666664
$stepSerial  = 0.0416666666666667

$bValues = @(
    44348.291666666664,
    44348.333333333336,
    44348.375000057873,
    44348.41666678241,
    44348.458333506947,
    44348.500000231485,
    44348.541666956022,
    44348.583333680559,
    44348.625000405096,
    44348.666667129626,
    44348.708333854163,
    44348.750000578701,
    44348.791667303238,
    44348.833334027775,
    44348.875000752312,
    44348.916667476849,
    44348.958334201387,
    44349.000000925924,
    44349.041667650461,
    44349.083334374998,
    44349.125001099535,
    44349.166667824073,
    44349.20833454861,
    44349.250001273147
)

$cValues = @(
    44348.333333333336,
    44348.375000057873,
    44348.41666678241,
    44348.458333506947,
    44348.500000231485,
    44348.541666956022,
    44348.583333680559,
    44348.625000405096,
    44348.666667129626,
    44348.708333854163,
    44348.750000578701,
    44348.791667303238,
    44348.833334027775,
    44348.875000752312,
    44348.916667476849,
    44348.958334201387,
    44349.000000925924,
    44349.041667650461,
    44349.083334374998,
    44349.125001099535,
    44349.166667824073,
    44349.20833454861,
    44349.250001273147,
    44349.291667939811
)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $wsTime.Cells.Item($row, 2).Value = $bValues[$i]
    $wsTime.Cells.Item($row, 3).Value = $cValues[$i]
}

$rngTime = $wsTime.Range("B2:C25")
$rngTime.ClearFormats()
$rngTime.NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# Columns B and C now share a single width
$wsTime.Columns.Item(2).ColumnWidth = 18.109375
$wsTime.Columns.Item(3).ColumnWidth = 18.109375

# Update selection on this sheet and drop its "tabSelected" state
$wsTime.Range("C23:C25").Select()

# ---------------------------------------------------------------------
# EVData sheet: highlight column D (MaxFastChargingPower) and select it
# ---------------------------------------------------------------------
$wsEv = $wb.Worksheets.Item("EVData")

$rngEv = $wsEv.Range("D1:D21")
$rngEv.Interior.Color = 65535

$wsEv.Activate()
$wsEv.Range("G12").Select()
